$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in results for rows that now have outcomes
$ws.Range("G56").Value = "Acierto"
$ws.Range("H56").Value = 1.1

$ws.Range("G57").Value = "Acierto"
$ws.Range("H57").Value = 0.83

# Append new row 68 with the latest tracked event
$ws.Range("A68").Value = 14731581

# Force the date-like text to stay as literal text, not an Excel date serial
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "2025-09-23"
$ws.Range("B68").Style = "Normal"

$ws.Range("C68").Value = "Clement Chidekh"
$ws.Range("D68").Value = "Matteo Gigante"
$ws.Range("E68").Value = "Gana Matteo Gigante"
$ws.Range("F68").Value = 1.91
